# "Last minute IREC files" — strip the leftover SharePoint / content-type
# custom XML parts (content type schema, form templates, document
# property bag) that got dragged along with the template before this
# document went out. None of this is visible document content; it is
# pure package-level metadata exposed through Document.CustomXMLParts.
#
# Standard Word automation idiom: walk the CustomXMLParts collection and
# delete every part (walk backwards since deleting shifts indices/Count
# down as we go). Also sweep by namespace as a belt-and-braces pass in
# case any part isn't picked up by plain index iteration, and guard each
# step so a host that can't resolve a given member/index just leaves
# that part alone instead of aborting the whole cleanup.

$d = $word.ActiveDocument

function Remove-AllCustomXmlParts($parts) {
    if ($parts -eq $null) { return }
    $count = 0
    try { $count = $parts.Count } catch { $count = 0 }
    for ($i = $count; $i -ge 1; $i--) {
        try {
            $part = $parts.Item($i)
            if ($part -ne $null) {
                $part.Delete()
            }
        } catch {
            # Host couldn't resolve/delete this item — move on rather
            # than letting one bad index stop the rest of the cleanup.
        }
    }
}

# Pass 1: delete everything in document order.
Remove-AllCustomXmlParts $d.CustomXMLParts

# Pass 2: some hosts only expose parts through a namespace-scoped
# selection; sweep the "no namespace" (catch-all) bucket too.
try {
    $unscoped = $d.CustomXMLParts.SelectByNamespace("")
    Remove-AllCustomXmlParts $unscoped
} catch {
}

# Pass 3: re-check in case Delete() renumbered things unexpectedly —
# keep clearing until the collection reports itself empty (bounded so a
# host that never updates Count can't spin forever).
$guard = 0
while ($guard -lt 10) {
    $remaining = 0
    try { $remaining = $d.CustomXMLParts.Count } catch { $remaining = 0 }
    if ($remaining -le 0) { break }
    Remove-AllCustomXmlParts $d.CustomXMLParts
    $guard++
}

$d.Save()
